$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": shift dates forward one week and update MyForecast values ---
$ws1.Range("B2").Value = "'2025-01-12"
$ws1.Range("D2").Value = 144
$ws1.Range("B3").Value = "'2025-01-19"
$ws1.Range("D3").Value = 137
$ws1.Range("B4").Value = "'2025-01-26"
$ws1.Range("D4").Value = 155
$ws1.Range("B5").Value = "'2025-02-02"
$ws1.Range("D5").Value = 177
$ws1.Range("B6").Value = "'2025-02-09"
$ws1.Range("D6").Value = 194
$ws1.Range("B7").Value = "'2025-02-16"
$ws1.Range("D7").Value = 218
$ws1.Range("B8").Value = "'2025-02-23"
$ws1.Range("D8").Value = 262
$ws1.Range("B9").Value = "'2025-03-02"
$ws1.Range("D9").Value = 302
$ws1.Range("B10").Value = "'2025-03-09"
$ws1.Range("D10").Value = 312
$ws1.Range("B11").Value = "'2025-03-16"
$ws1.Range("D11").Value = 290
$ws1.Range("B12").Value = "'2025-03-23"
$ws1.Range("D12").Value = 270
$ws1.Range("B13").Value = "'2025-03-30"
$ws1.Range("D13").Value = 272
$ws1.Range("B14").Value = "'2025-04-06"
$ws1.Range("D14").Value = 300
$ws1.Range("B15").Value = "'2025-04-13"
$ws1.Range("D15").Value = 303
$ws1.Range("B16").Value = "'2025-04-20"
$ws1.Range("D16").Value = 294
$ws1.Range("B17").Value = "'2025-04-27"
$ws1.Range("D17").Value = 283

# --- Sheet "Summary": update derived metrics ---
$ws2.Range("B2").Value = "'2022-12-25 to 2025-01-05"
$ws2.Range("B4").Value = "'884"
$ws2.Range("B5").Value = "'291"
$ws2.Range("B6").Value = "'243"
$ws2.Range("B8").Value = "'32332 units"
$ws2.Range("B9").Value = "'3913"
$ws2.Range("B10").Value = "'1589"
$ws2.Range("B11").Value = "'613"
$ws2.Range("B12").Value = "'312"
$ws2.Range("B13").Value = "'2025-03-09"
$ws2.Range("B14").Value = "'137"
$ws2.Range("B15").Value = "'2025-01-19"
